$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 35657
$ws.Range("J3").Value = 35657
$ws.Range("L3").Value = 35657
$ws.Range("N3").Value = -35885
$ws.Range("H102").Value = 35657
$ws.Range("J102").Value = 35657
$ws.Range("L102").Value = 35657
$ws.Range("N102").Value = -42147
$ws.Range("H121").Value = 1025
$ws.Range("J121").Value = 1900
$ws.Range("L121").Value = 5700
$ws.Range("N121").Value = -9194
$ws.Range("H129").Value = 4665.346
$ws.Range("J129").Value = 1028.25
$ws.Range("L129").Value = 3084.75
$ws.Range("N129").Value = -13084.75
$ws.Range("H137").Value = 1938.409
$ws.Range("I137").Value = 1403
$ws.Range("J137").Value = 3085.7144
$ws.Range("K137").Value = 4209
$ws.Range("L137").Value = 9257.143199999999
$ws.Range("M137").Value = -1659
$ws.Range("N137").Value = -14357.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28288.11
$ws.Range("I32").Value = 5280.2183
$ws.Range("J32").Value = 168891.89
$ws.Range("K32").Value = 5280.2183
$ws.Range("L32").Value = 168891.89
$ws.Range("M32").Value = -4993.2183
$ws.Range("N32").Value = -169465.89
$ws.Range("H61").Value = 1570.8182
$ws.Range("I61").Value = 991.8333
$ws.Range("J61").Value = 2265.6
$ws.Range("K61").Value = 991.8333
$ws.Range("L61").Value = 2265.6
$ws.Range("M61").Value = -779.8333
$ws.Range("N61").Value = -2689.6
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H136").Value = 1570.8182
$ws.Range("I136").Value = 991.8333
$ws.Range("J136").Value = 2265.6
$ws.Range("K136").Value = 2975.4999
$ws.Range("L136").Value = 6796.799999999999
$ws.Range("M136").Value = -425.4998999999998
$ws.Range("N136").Value = -11896.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 75067.47
$ws.Range("I86").Value = 86216.30499999999
$ws.Range("J86").Value = 2600
$ws.Range("K86").Value = 86216.30499999999
$ws.Range("L86").Value = 2600
$ws.Range("M86").Value = -85093.30499999999
$ws.Range("N86").Value = -4846
$ws.Range("H89").Value = 75067.47
$ws.Range("I89").Value = 86216.30499999999
$ws.Range("J89").Value = 2600
$ws.Range("K89").Value = 431081.525
$ws.Range("L89").Value = 13000
$ws.Range("M89").Value = -425465.525
$ws.Range("N89").Value = -24232
$ws.Range("H134").Value = 2420.1191
$ws.Range("I134").Value = 2427.7878
$ws.Range("J134").Value = 2392
$ws.Range("K134").Value = 7283.3634
$ws.Range("L134").Value = 7176
$ws.Range("M134").Value = -4748.3634
$ws.Range("N134").Value = -12246

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 995.8461
$ws.Range("I22").Value = 461.5
$ws.Range("J22").Value = 1850.8
$ws.Range("K22").Value = 461.5
$ws.Range("L22").Value = 1850.8
$ws.Range("M22").Value = -111.5
$ws.Range("N22").Value = -2550.8
$ws.Range("H31").Value = 22887.055
$ws.Range("I31").Value = 894.88635
$ws.Range("J31").Value = 56254.484
$ws.Range("K31").Value = 894.88635
$ws.Range("L31").Value = 56254.484
$ws.Range("M31").Value = -599.88635
$ws.Range("N31").Value = -56844.484
$ws.Range("H34").Value = 22887.055
$ws.Range("I34").Value = 894.88635
$ws.Range("J34").Value = 56254.484
$ws.Range("K34").Value = 894.88635
$ws.Range("L34").Value = 56254.484
$ws.Range("M34").Value = -692.88635
$ws.Range("N34").Value = -56658.484
$ws.Range("H35").Value = 2268.182
$ws.Range("J35").Value = 4850
$ws.Range("L35").Value = 4850
$ws.Range("N35").Value = -5438
$ws.Range("H107").Value = 1363.4445
$ws.Range("I107").Value = 1647.1666
$ws.Range("J107").Value = 796
$ws.Range("K107").Value = 1647.1666
$ws.Range("L107").Value = 796
$ws.Range("M107").Value = 272.8334
$ws.Range("N107").Value = -4636
$ws.Range("H141").Value = 71025
$ws.Range("J141").Value = 44700
$ws.Range("L141").Value = 44700
$ws.Range("N141").Value = -55060

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 117.666664
$ws.Range("I38").Value = 9.166667
$ws.Range("J38").Value = 226.16667
$ws.Range("K38").Value = 27.500001
$ws.Range("L38").Value = 678.50001
$ws.Range("M38").Value = 319.499999
$ws.Range("N38").Value = -1372.50001
$ws.Range("H63").Value = 1307
$ws.Range("I63").Value = 1200
$ws.Range("J63").Value = 1414
$ws.Range("K63").Value = 3600
$ws.Range("L63").Value = 4242
$ws.Range("M63").Value = -2851
$ws.Range("N63").Value = -5740
$ws.Range("H66").Value = 1307
$ws.Range("I66").Value = 1200
$ws.Range("J66").Value = 1414
$ws.Range("K66").Value = 10800
$ws.Range("L66").Value = 12726
$ws.Range("M66").Value = -7056
$ws.Range("N66").Value = -20214
$ws.Range("H137").Value = 17377084
$ws.Range("J137").Value = 25492820
$ws.Range("L137").Value = 76478460
$ws.Range("N137").Value = -76488660
$ws.Range("H140").Value = 6173.826
$ws.Range("I140").Value = 8707
$ws.Range("K140").Value = 26121
$ws.Range("M140").Value = -20941

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 76629.21000000001
$ws.Range("I70").Value = 129012.44
$ws.Range("J70").Value = 6784.9165
$ws.Range("K70").Value = 129012.44
$ws.Range("L70").Value = 6784.9165
$ws.Range("M70").Value = -128742.44
$ws.Range("N70").Value = -7324.9165
$ws.Range("H73").Value = 76629.21000000001
$ws.Range("I73").Value = 129012.44
$ws.Range("J73").Value = 6784.9165
$ws.Range("K73").Value = 129012.44
$ws.Range("L73").Value = 6784.9165
$ws.Range("M73").Value = -128076.44
$ws.Range("N73").Value = -8656.916499999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2254.818
$ws.Range("I46").Value = 200.14285
$ws.Range("J46").Value = 5850.5
$ws.Range("K46").Value = 200.14285
$ws.Range("L46").Value = 5850.5
$ws.Range("M46").Value = -12.14285000000001
$ws.Range("N46").Value = -6226.5
$ws.Range("H103").Value = 31620.334
$ws.Range("J103").Value = 31620.334
$ws.Range("L103").Value = 31620.334
$ws.Range("N103").Value = -33964.334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 48249.5
$ws.Range("J135").Value = 48249.5
$ws.Range("L135").Value = 48249.5
$ws.Range("N135").Value = -58389.5
